$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (rows 2-6) from 10.6 to 3.65
$ws.Range("A2:A6").Value = 3.65
